$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume number, date range covered) ---
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -80

# --- Row 15 (Rape) ---
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -45.454545454545

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 83
$ws.Range("K16").Value = -1.204819277108
$ws.Range("L16").Value = 13.888888888888
$ws.Range("M16").Value = -23.364485981308
$ws.Range("N16").Value = -57.948717948717

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -30.434782608695
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = 17.592592592592
$ws.Range("L17").Value = 42.696629213483
$ws.Range("M17").Value = 56.790123456790
$ws.Range("N17").Value = 29.591836734693

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 120
$ws.Range("I18").Value = 67
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 48.888888888888
$ws.Range("L18").Value = 59.523809523809
$ws.Range("M18").Value = -51.449275362318
$ws.Range("N18").Value = -82.597402597402

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 32.558139534883
$ws.Range("I19").Value = 298
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 24.166666666666
$ws.Range("L19").Value = 69.318181818181
$ws.Range("M19").Value = 64.640883977900
$ws.Range("N19").Value = 58.510638297872

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 275
$ws.Range("F20").Value = 52
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 147.619047619048
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 113
$ws.Range("K20").Value = 120.353982300885
$ws.Range("L20").Value = 167.741935483871
$ws.Range("M20").Value = 227.631578947368
$ws.Range("N20").Value = -70.462633451957

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 119
$ws.Range("H21").Value = 30.252100840336
$ws.Range("I21").Value = 830
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 38.333333333333
$ws.Range("L21").Value = 72.557172557172
$ws.Range("M21").Value = 40.202702702702
$ws.Range("N21").Value = -51.884057971014

# --- Row 22 (Transit) --- C22 switches from "N/A" text to a real number
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 12.5
$ws.Range("L22").Value = 125
$ws.Range("M22").Value = 12.5

# --- Row 23 (Housing) --- D23/E23 switch from numbers to "N/A" (text) placeholders
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 400
$ws.Range("I23").Value = 27
$ws.Range("K23").Value = 92.857142857142
$ws.Range("L23").Value = 50
$ws.Range("M23").Value = 58.823529411764

# --- Row 24 (Petit Larceny) ---
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -19.354838709677
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = 4.901960784313
$ws.Range("I24").Value = 605
$ws.Range("J24").Value = 516
$ws.Range("K24").Value = 17.248062015503
$ws.Range("L24").Value = 51.25
$ws.Range("M24").Value = 4.671280276816

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 34
$ws.Range("H25").Value = 17.241379310344
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 198
$ws.Range("K25").Value = 11.111111111111
$ws.Range("L25").Value = 20.218579234972
$ws.Range("M25").Value = 22.222222222222

# --- Row 26 (UCR Rape*) ---
$ws.Range("D26").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = -31.578947368421
$ws.Range("L26").Value = 85.714285714285

# --- Row 27 (Other Sex Crimes) --- C27 switches from "N/A" text to a real number
$ws.Range("C27").Value = 2
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -24
$ws.Range("L27").Value = -5

# --- Row 28 (Shooting Vic.) ---
$ws.Range("N28").Value = -75

# --- Row 29 (Shooting Inc.) ---
$ws.Range("N29").Value = -75

# --- Row 30 (Hate Crimes) --- D30/E30 switch from numbers to "N/A" (text) placeholders
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
